$wb = $excel.ActiveWorkbook

# Duplicate the "Portugal" worksheet (closest template for a new market sheet)
# and place the copy immediately after it, then rename to "Slovakia".
$port = $wb.Worksheets.Item("Portugal")
$port.Activate()
$port.Range("A1:XFD1048576").Select()

$port.Copy([System.Reflection.Missing]::Value, $port)
$newSheet = $wb.Worksheets.Item($port.Index + 1)
$newSheet.Name = "Slovakia"

# Update market-specific content
$newSheet.Range("B2").Value = "Slovakia Market"
$newSheet.Range("B4").Value = "NGC-2930/T3177"

# Rows 3 & 4 should use the default row height on the new sheet
$newSheet.Rows.Item(3).AutoFit()
$newSheet.Rows.Item(4).AutoFit()

# Set the view/selection state of the new sheet
$newSheet.Range("E14").Select()
